$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("Objetivos:"): fill in the responsible professor for the section
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# Row 13 ("Programa resumido:"): fill with the activation date text
# (force text, not an auto-converted date, then restore the original
#  cell formatting so the style index used by the cell stays unchanged)
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "01/01/2023"
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "01/01/2023"
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# Row 15 ("Programa:"): fill with the other responsible professor
$ws.Range("B15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value = "519033 - Carlos Yujiro Shigue"

# Row 18 ("Método:"): fill with the second responsible professor
$ws.Range("B18").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C18").Value = "7290967 - Emerson Gonçalves de Melo"
